$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J12").Value = 149.8
$ws.Range("K12").Value = 1
$ws.Range("M12").Value = 0
$ws.Range("J14").Value = 82
$ws.Range("K14").Value = 1
$ws.Range("M14").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("AC15").Value = 0
$ws.Range("J16").Value = 47.9
$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 0
$ws.Range("J17").Value = 47.5
$ws.Range("K17").Value = 1
$ws.Range("M17").Value = 0
$ws.Range("E25").Value = 4
$ws.Range("I25").Value = 3
$ws.Range("M25").Value = 1
$ws.Range("Q25").Value = 3
$ws.Range("U25").Value = 4
$ws.Range("R28").Value = 58.8
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 2
$ws.Range("U31").Value = 0
$ws.Range("U32").Value = 11
$ws.Range("R35").Value = 438.0786219999993
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("U36").Value = 0
$ws.Range("J37").Value = 407.320896
$ws.Range("B40").Value = 197.8405019999992
$ws.Range("F40").Value = 424.6197320946533
$ws.Range("N40").Value = 268.5294299999971
$ws.Range("V40").Value = 432.2194689189191
$ws.Range("J41").Value = 54.1800000000014
$ws.Range("E42").Value = 4
$ws.Range("Q42").Value = 5
$ws.Range("F46").Value = 79.96320590534793
$ws.Range("M46").Value = 11
$ws.Range("Q46").Value = 0
$ws.Range("V46").Value = 87.36028308108128
$ws.Range("M47").Value = 5
$ws.Range("U47").Value = 6
$ws.Range("N50").Value = 47.28000000000103
$ws.Range("O50").Value = 1
$ws.Range("Q50").Value = 0
$ws.Range("E51").Value = 10
$ws.Range("J51").Value = 30.7
$ws.Range("K51").Value = 1
$ws.Range("L51").Value = 2
$ws.Range("N51").Value = 30.7
$ws.Range("O51").Value = 1
$ws.Range("P51").Value = 12
$ws.Range("X51").Value = 3
$ws.Range("Z51").Value = 30.9400000000002
$ws.Range("J52").Value = 29.8
$ws.Range("K52").Value = 1
$ws.Range("M52").Value = 0
$ws.Range("R52").Value = 0
$ws.Range("S52").Value = 0
$ws.Range("U52").Value = 5
$ws.Range("M53").Value = 4
$ws.Range("Q53").Value = 5
$ws.Range("I54").Value = 3
$ws.Range("M54").Value = 3
$ws.Range("Q54").Value = 4
$ws.Range("I56").Value = 2
$ws.Range("M56").Value = 1
$ws.Range("Q56").Value = 3
$ws.Range("I58").Value = 5
$ws.Range("M58").Value = 3
$ws.Range("Q58").Value = 5
